$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 16:47:39"
$ws.Range("K2").Value = "8.7 MJ/m2"
$ws.Range("E3").Value = "2026-02-06 16:47:42"
$ws.Range("K3").Value = "12.4 MJ/m2"
$ws.Range("E4").Value = "2026-02-06 16:47:44"
$ws.Range("J4").Value = "996.9 hPa"
$ws.Range("K4").Value = "11.4 MJ/m2"
$ws.Range("E5").Value = "2026-02-06 16:47:47"
$ws.Range("J5").Value = "997.2 hPa"
$ws.Range("O5").Value = "11.0 °C"
$ws.Range("E6").Value = "2026-02-06 16:47:49"
$ws.Range("J6").Value = "998.4 hPa"
$ws.Range("K6").Value = "10.0 MJ/m2"
$ws.Range("E7").Value = "2026-02-06 16:47:51"
$ws.Range("J7").Value = "998.0 hPa"
$ws.Range("K7").Value = "12.2 MJ/m2"
$ws.Range("O7").Value = "11.7 °C"
$ws.Range("E8").Value = "2026-02-06 16:47:54"
$ws.Range("K8").Value = "11.8 MJ/m2"
$ws.Range("O8").Value = "10.2 °C"
$ws.Range("E9").Value = "2026-02-06 16:47:56"
$ws.Range("H9").Value = "'84%"
$ws.Range("O9").Value = "4.9 °C"
$ws.Range("E10").Value = "2026-02-06 16:47:58"
$ws.Range("O10").Value = "9.1 °C"
$ws.Range("E11").Value = "2026-02-06 16:48:01"
$ws.Range("H11").Value = "'76%"
$ws.Range("J11").Value = "998.4 hPa"
$ws.Range("K11").Value = "9.0 MJ/m2"
$ws.Range("O11").Value = "5.4 °C"
$ws.Range("E12").Value = "2026-02-06 16:48:03"
$ws.Range("H12").Value = "'56%"
$ws.Range("K12").Value = "12.3 MJ/m2"
$ws.Range("E13").Value = "2026-02-06 16:48:06"
$ws.Range("O13").Value = "10.3 °C"
$ws.Range("E14").Value = "2026-02-06 16:48:08"
$ws.Range("K14").Value = "7.5 MJ/m2"
$ws.Range("E15").Value = "2026-02-06 16:48:10"
$ws.Range("H15").Value = "'71%"
$ws.Range("J15").Value = "997.3 hPa"
$ws.Range("K15").Value = "11.6 MJ/m2"
$ws.Range("O15").Value = "10.6 °C"
$ws.Range("E16").Value = "2026-02-06 16:48:13"
$ws.Range("K16").Value = "9.7 MJ/m2"
$ws.Range("O16").Value = "6.0 °C"
$ws.Range("E17").Value = "2026-02-06 16:48:15"
$ws.Range("H17").Value = "'84%"
$ws.Range("K17").Value = "10.6 MJ/m2"
$ws.Range("O17").Value = "6.0 °C"
$ws.Range("E18").Value = "2026-02-06 16:48:17"
$ws.Range("K18").Value = "5.9 MJ/m2"
$ws.Range("E19").Value = "2026-02-06 16:48:20"
$ws.Range("H19").Value = "'76%"
$ws.Range("J19").Value = "999.6 hPa"
$ws.Range("K19").Value = "11.7 MJ/m2"
$ws.Range("O19").Value = "9.9 °C"
$ws.Range("E20").Value = "2026-02-06 16:48:23"
$ws.Range("H20").Value = "'78%"
$ws.Range("K20").Value = "12.2 MJ/m2"
$ws.Range("E21").Value = "2026-02-06 16:48:25"
$ws.Range("H21").Value = "'73%"
$ws.Range("J21").Value = "997.5 hPa"
$ws.Range("K21").Value = "10.3 MJ/m2"
$ws.Range("O21").Value = "8.5 °C"
$ws.Range("E22").Value = "2026-02-06 16:48:27"
$ws.Range("O22").Value = "10.8 °C"
$ws.Range("E23").Value = "2026-02-06 16:48:30"
$ws.Range("J23").Value = "997.3 hPa"
$ws.Range("O23").Value = "10.1 °C"
$ws.Range("E24").Value = "2026-02-06 16:48:32"
$ws.Range("J24").Value = "996.7 hPa"
$ws.Range("K24").Value = "11.8 MJ/m2"
$ws.Range("E25").Value = "2026-02-06 16:48:35"
$ws.Range("H25").Value = "'78%"
$ws.Range("K25").Value = "9.8 MJ/m2"
$ws.Range("O25").Value = "4.4 °C"
$ws.Range("E26").Value = "2026-02-06 16:48:37"
$ws.Range("I26").Value = "0.3 mm"
$ws.Range("K26").Value = "8.0 MJ/m2"
$ws.Range("O26").Value = "-0.8 °C"
$ws.Range("E27").Value = "2026-02-06 16:48:40"
$ws.Range("H27").Value = "'80%"
$ws.Range("J27").Value = "997.3 hPa"
$ws.Range("K27").Value = "10.6 MJ/m2"
$ws.Range("O27").Value = "10.9 °C"
$ws.Range("E28").Value = "2026-02-06 16:48:42"
$ws.Range("O28").Value = "4.8 °C"
$ws.Range("E29").Value = "2026-02-06 16:48:44"
$ws.Range("K29").Value = "12.1 MJ/m2"
$ws.Range("O29").Value = "12.8 °C"
$ws.Range("E30").Value = "2026-02-06 16:48:47"
$ws.Range("H30").Value = "'76%"
$ws.Range("K30").Value = "9.0 MJ/m2"
$ws.Range("L30").Value = "43.9 km/h - 189º 16:13 TU"
$ws.Range("O30").Value = "-3.6 °C"
$ws.Range("E31").Value = "2026-02-06 16:48:49"
$ws.Range("O31").Value = "7.3 °C"
$ws.Range("E32").Value = "2026-02-06 16:48:52"
$ws.Range("H32").Value = "'49%"
$ws.Range("J32").Value = "998.7 hPa"
$ws.Range("O32").Value = "15.8 °C"
$ws.Range("E33").Value = "2026-02-06 16:48:54"
$ws.Range("O33").Value = "10.2 °C"
$ws.Range("E34").Value = "2026-02-06 16:48:56"
$ws.Range("K34").Value = "12.0 MJ/m2"
$ws.Range("O34").Value = "8.7 °C"
$ws.Range("E35").Value = "2026-02-06 16:48:59"
$ws.Range("K35").Value = "9.5 MJ/m2"
$ws.Range("E36").Value = "2026-02-06 16:49:01"
$ws.Range("I36").Value = "0.9 mm"
$ws.Range("J36").Value = "999.6 hPa"
$ws.Range("L36").Value = "36.7 km/h - 182º 16:06 TU"
